$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Rows 66 / 67: the two fixtures' F:V (match) data were swapped while A:E
# (Indice/pais/torneio/temporada/data_partida) stay put on their own row.
# ---------------------------------------------------------------------------
$ws.Cells.Item(66,6).Value  = "Cardiff Metropolitan"
$ws.Cells.Item(66,7).Value  = 2
$ws.Cells.Item(66,8).Value  = "Newtown"
$ws.Cells.Item(66,9).Value  = 1
$ws.Cells.Item(66,10).Value = 3.25
$ws.Cells.Item(66,11).Value = "16/10/2023 08:12"
$ws.Cells.Item(66,12).Value = 3.15
$ws.Cells.Item(66,13).Value = "17/10/2023 20:39"
$ws.Cells.Item(66,14).Value = 3.36
$ws.Cells.Item(66,15).Value = "16/10/2023 08:12"
$ws.Cells.Item(66,16).Value = 3.24
$ws.Cells.Item(66,17).Value = "17/10/2023 20:39"
$ws.Cells.Item(66,18).Value = 2.06
$ws.Cells.Item(66,19).Value = "16/10/2023 08:12"
$ws.Cells.Item(66,20).Value = 2.31
$ws.Cells.Item(66,21).Value = "17/10/2023 20:39"
$ws.Cells.Item(66,22).Value = "https://www.betexplorer.com/football/wales/cymru-premier/cardiff-metropolitan-university-newtown/Q5G6CKXH/"

$ws.Cells.Item(67,6).Value  = "Haverfordwest"
$ws.Cells.Item(67,7).Value  = 3
$ws.Cells.Item(67,8).Value  = "Penybont"
$ws.Cells.Item(67,9).Value  = 2
$ws.Cells.Item(67,10).Value = 3.39
$ws.Cells.Item(67,11).Value = "16/10/2023 16:42"
$ws.Cells.Item(67,12).Value = 3.12
$ws.Cells.Item(67,13).Value = "17/10/2023 20:24"
$ws.Cells.Item(67,14).Value = 3.36
$ws.Cells.Item(67,15).Value = "16/10/2023 16:42"
$ws.Cells.Item(67,16).Value = 3.39
$ws.Cells.Item(67,17).Value = "17/10/2023 18:49"
$ws.Cells.Item(67,18).Value = 1.96
$ws.Cells.Item(67,19).Value = "16/10/2023 16:42"
$ws.Cells.Item(67,20).Value = 2.26
$ws.Cells.Item(67,21).Value = "17/10/2023 20:24"
$ws.Cells.Item(67,22).Value = "https://www.betexplorer.com/football/wales/cymru-premier/haverfordwest-penybont/K6YvfCHj/"

# ---------------------------------------------------------------------------
# Rows 75 / 76: same kind of swap.
# ---------------------------------------------------------------------------
$ws.Cells.Item(75,6).Value  = "Penybont"
$ws.Cells.Item(75,7).Value  = 0
$ws.Cells.Item(75,8).Value  = "Aberystwyth"
$ws.Cells.Item(75,9).Value  = 2
$ws.Cells.Item(75,10).Value = 1.26
$ws.Cells.Item(75,12).Value = 1.46
$ws.Cells.Item(75,14).Value = 5.37
$ws.Cells.Item(75,16).Value = 4.44
$ws.Cells.Item(75,17).Value = "27/10/2023 20:36"
$ws.Cells.Item(75,18).Value = 8.199999999999999
$ws.Cells.Item(75,20).Value = 6.86
$ws.Cells.Item(75,21).Value = "27/10/2023 20:36"
$ws.Cells.Item(75,22).Value = "https://www.betexplorer.com/football/wales/cymru-premier/penybont-aberystwyth/EsDH7Voj/"

$ws.Cells.Item(76,6).Value  = "Connahs Q."
$ws.Cells.Item(76,7).Value  = 6
$ws.Cells.Item(76,8).Value  = "Caernarfon"
$ws.Cells.Item(76,9).Value  = 1
$ws.Cells.Item(76,10).Value = 1.36
$ws.Cells.Item(76,12).Value = 1.38
$ws.Cells.Item(76,14).Value = 4.78
$ws.Cells.Item(76,16).Value = 5.04
$ws.Cells.Item(76,17).Value = "27/10/2023 20:41"
$ws.Cells.Item(76,18).Value = 6.43
$ws.Cells.Item(76,20).Value = 7.64
$ws.Cells.Item(76,21).Value = "27/10/2023 20:41"
$ws.Cells.Item(76,22).Value = "https://www.betexplorer.com/football/wales/cymru-premier/connahs-q-caernarfon/KKDD8BWq/"

# ---------------------------------------------------------------------------
# New row 81: next fixture appended at the bottom of the sheet.
# ---------------------------------------------------------------------------
$ws.Cells.Item(81,1).Value  = 80
$ws.Cells.Item(81,2).Value  = "wales"
$ws.Cells.Item(81,3).Value  = "cymru-premier"
$ws.Cells.Item(81,4).Value  = "2023-2024"
$ws.Cells.Item(81,5).Value  = 45230.86458333334
$ws.Cells.Item(81,6).Value  = "TNS"
$ws.Cells.Item(81,7).Value  = 6
$ws.Cells.Item(81,8).Value  = "Colwyn Bay"
$ws.Cells.Item(81,9).Value  = 1
$ws.Cells.Item(81,10).Value = 1.03
$ws.Cells.Item(81,11).Value = "30/10/2023 16:42"
$ws.Cells.Item(81,12).Value = 1.02
$ws.Cells.Item(81,13).Value = "31/10/2023 19:58"
$ws.Cells.Item(81,14).Value = 16.84
$ws.Cells.Item(81,15).Value = "30/10/2023 16:42"
$ws.Cells.Item(81,16).Value = 27.01
$ws.Cells.Item(81,17).Value = "31/10/2023 20:29"
$ws.Cells.Item(81,18).Value = 22.13
$ws.Cells.Item(81,19).Value = "30/10/2023 16:42"
$ws.Cells.Item(81,20).Value = 47.31
$ws.Cells.Item(81,21).Value = "31/10/2023 20:29"
$ws.Cells.Item(81,22).Value = "https://www.betexplorer.com/football/wales/cymru-premier/tns-colwyn-bay/Y9MMIefi/"

# Mirror the formatting of the row above (border/bold index cell + date style)
# onto the freshly added row, the same way it is applied to every other row.
$ws.Cells.Item(80,1).Copy()
$ws.Cells.Item(81,1).PasteSpecial(-4122)

$ws.Cells.Item(80,5).Copy()
$ws.Cells.Item(81,5).PasteSpecial(-4122)

$excel.CutCopyMode = 0
